$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the row above (A4) into the new row's label cell (A5),
# then set its text value.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "2021年"

# Fill in the numeric data for the new row.
$ws.Range("B5").Value = 7.5
$ws.Range("C5").Value = -10.9
$ws.Range("D5").Value = -3.4
$ws.Range("E5").Value = -3.8
$ws.Range("F5").Value = 4.4
$ws.Range("G5").Value = 5.8
